$d = $word.ActiveDocument

# Build the new contact-info paragraph text next to a genuinely "blank"
# paragraph (the empty one under CORE COMPETENCIES) so the freshly
# inserted paragraph/run picks up no inherited direct formatting
# (no bold/size run properties, no paragraph style).
$blankPara = $d.Paragraphs(5)
$blankPara.Range.InsertParagraphBefore()
$stagingPara = $d.Paragraphs(5)
$stagingPara.Range.Text = "202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX"
$stagingPara.Alignment = 1

# Move that cleanly-formatted paragraph (text + its paragraph mark) to sit
# immediately after the "Dheeraj Chand" title paragraph.
$stagingPara.Range.Cut()

$titlePara = $d.Paragraphs.First
$insertPoint = $d.Range($titlePara.Range.End, $titlePara.Range.End)
$insertPoint.Paste()
